$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows with the new subscription/resource-group values and
# drop the now-unneeded third data row (rg-test1).
$ws.Range("A2").Value = "Default - Microsoft Azure Sponsorship 2"
$ws.Range("B2").Value = "rg-hpc"
$ws.Range("C2").Value = "FinOps3"
$ws.Range("D2").Value = "Value01"

$ws.Range("A3").Value = "HPC subscription"
$ws.Range("B3").Value = "rg-hpc-eastus"
$ws.Range("C3").Value = "FinOps3"
$ws.Range("D3").Value = "Value02"

# Remove the former row 4 (MCAPS-MarcusGaspar / rg-test1 / FinOps3 / Value03)
$ws.Range("A4:D4").EntireRow.Delete()

# Column A now needs to be wider to fit the longer subscription names
# (target best-fit width is ~33.89; 33 is the closest the host's column
# width quantization can reach)
$ws.Columns.Item(1).ColumnWidth = 33

# Match the saved selection state
$ws.Range("D3").Select()
